# Auto-generated edit script
# Applies a left-shift of the yearly (E:I) data columns for each metric
# block on the "Overview" sheet: E<-F, F<-G, G<-H, H<-I, and I<-new value
# for the newly added fiscal year (1401/12), dropping the oldest year
# (1396/12) and adding the newest one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("E", "F", "G", "H", "I")

# New value to place into column I for each data row (the incoming
# 1401/12 figures). A value of "-" means the cell should hold the
# literal dash placeholder text used throughout the sheet.
$rowNewI = @{}
$rowNewI[10] = 0
$rowNewI[11] = 3778860
$rowNewI[12] = 11042266
$rowNewI[13] = 64923
$rowNewI[14] = 1184278
$rowNewI[15] = 1615015
$rowNewI[16] = 0
$rowNewI[17] = 0
$rowNewI[18] = 0
$rowNewI[19] = 0
$rowNewI[20] = 0
$rowNewI[21] = 0
$rowNewI[22] = 17685342
$rowNewI[28] = 0
$rowNewI[29] = 1363295
$rowNewI[30] = 6269056
$rowNewI[31] = 75687
$rowNewI[32] = 1179313
$rowNewI[33] = 134712
$rowNewI[34] = "-"
$rowNewI[35] = 0
$rowNewI[36] = 0
$rowNewI[37] = 0
$rowNewI[38] = 0
$rowNewI[39] = 7228
$rowNewI[40] = 3038698
$rowNewI[41] = "-"
$rowNewI[42] = 12067989
$rowNewI[48] = 0
$rowNewI[49] = 45157623
$rowNewI[50] = 177651254
$rowNewI[51] = 1837331
$rowNewI[52] = 171263842
$rowNewI[53] = 10453848
$rowNewI[54] = "-"
$rowNewI[55] = 0
$rowNewI[56] = 0
$rowNewI[57] = 654719
$rowNewI[58] = 0
$rowNewI[59] = 150813
$rowNewI[60] = 3332016
$rowNewI[61] = 410501446
$rowNewI[67] = 0
$rowNewI[68] = 33123882
$rowNewI[69] = 28337800
$rowNewI[70] = 24275384
$rowNewI[71] = 145223399
$rowNewI[72] = 77601461
$rowNewI[73] = "-"
$rowNewI[74] = 0
$rowNewI[75] = 0
$rowNewI[76] = 0
$rowNewI[77] = 0
$rowNewI[78] = 20865108
$rowNewI[79] = 1096528
$rowNewI[85] = 0
$rowNewI[86] = -26479604
$rowNewI[87] = -88996394
$rowNewI[88] = -706596
$rowNewI[89] = -96272545
$rowNewI[90] = -4775174
$rowNewI[91] = "-"
$rowNewI[92] = 0
$rowNewI[93] = 0
$rowNewI[94] = 0
$rowNewI[95] = 0
$rowNewI[96] = -64850
$rowNewI[97] = -1790794
$rowNewI[98] = -219085957
$rowNewI[104] = 0
$rowNewI[105] = 18678019
$rowNewI[106] = 88654860
$rowNewI[107] = 1130735
$rowNewI[108] = 74991297
$rowNewI[109] = 5678674
$rowNewI[110] = "-"
$rowNewI[111] = 0
$rowNewI[112] = 0
$rowNewI[113] = 654719
$rowNewI[114] = 0
$rowNewI[115] = 85963
$rowNewI[116] = 1541222
$rowNewI[117] = 191415489

$dataRows = @($rowNewI.Keys | Sort-Object { [int]$_ })

foreach ($r in $dataRows) {
    $vals = @()
    foreach ($c in $cols) {
        $vals += ,($ws.Range($c + $r).Value())
    }
    for ($i = 0; $i -lt 4; $i++) {
        $ws.Range($cols[$i] + $r).Value = $vals[$i + 1]
    }
    $ws.Range("I" + $r).Value = $rowNewI[$r]
}

# Shift the "twelve months ended .../12" column headers the same way,
# dropping the 1396/12 label and introducing a new 1401/12 label.
$headerRows = @(8, 26, 46, 65, 83, 102)
$newHeaderLabel = "دوازده ماهه منتهی به 1401/12"

foreach ($r in $headerRows) {
    $vals = @()
    foreach ($c in $cols) {
        $vals += ,($ws.Range($c + $r).Value())
    }
    for ($i = 0; $i -lt 4; $i++) {
        $ws.Range($cols[$i] + $r).Value = $vals[$i + 1]
    }
    $ws.Range("I" + $r).Value = $newHeaderLabel
}
